$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 169.28572
$ws.Range("I5").Value = 182.5
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 182.5
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = -67.5
$ws.Range("N5").Value = -320
# Row 15
$ws.Range("H15").Value = 2433.8276
$ws.Range("I15").Value = 2433.8276
$ws.Range("K15").Value = 7301.4828
$ws.Range("M15").Value = -7132.4828
# Row 40
$ws.Range("H40").Value = 4694.7334
$ws.Range("I40").Value = 5499.6665
$ws.Range("J40").Value = 4605.2964
$ws.Range("K40").Value = 5499.6665
$ws.Range("L40").Value = 4605.2964
$ws.Range("M40").Value = -5324.6665
$ws.Range("N40").Value = -4955.2964
# Row 55
$ws.Range("H55").Value = 468.2
$ws.Range("I55").Value = 826.3333
$ws.Range("J55").Value = 314.7143
$ws.Range("K55").Value = 826.3333
$ws.Range("L55").Value = 314.7143
$ws.Range("M55").Value = -612.3333
$ws.Range("N55").Value = -742.7143
# Row 58
$ws.Range("H58").Value = 441.5
$ws.Range("I58").Value = 441.5
$ws.Range("K58").Value = 1324.5
$ws.Range("M58").Value = -1174.5
# Row 103
$ws.Range("H103").Value = 5103525.5
$ws.Range("I103").Value = 2266.625
$ws.Range("J103").Value = 11905204
$ws.Range("K103").Value = 6799.875
$ws.Range("L103").Value = 35715612
$ws.Range("M103").Value = -6213.875
$ws.Range("N103").Value = -35716784
# Row 106
$ws.Range("H106").Value = 3500
$ws.Range("I106").Value = 3500
$ws.Range("K106").Value = 3500
$ws.Range("M106").Value = -2869
# Row 107
$ws.Range("H107").Value = 310.7857
$ws.Range("I107").Value = 310.7857
$ws.Range("K107").Value = 310.7857
$ws.Range("M107").Value = 1609.2143
# Row 113
$ws.Range("H113").Value = 5378
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -1746
# Row 132
$ws.Range("H132").Value = 4735.5757
$ws.Range("I132").Value = 4852.3125
$ws.Range("K132").Value = 14556.9375
$ws.Range("M132").Value = -12026.9375
# Row 137
$ws.Range("H137").Value = 1394286.8
$ws.Range("I137").Value = 2001903
$ws.Range("J137").Value = 13340.728
$ws.Range("K137").Value = 6005709
$ws.Range("L137").Value = 40022.18399999999
$ws.Range("M137").Value = -6003159
$ws.Range("N137").Value = -45122.18399999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
# Row 32
$ws.Range("H32").Value = 3857.3845
$ws.Range("I32").Value = 3891.76
$ws.Range("K32").Value = 3891.76
$ws.Range("M32").Value = -3604.76
# Row 43
$ws.Range("H43").Value = 32742.25
$ws.Range("J43").Value = 32742.25
$ws.Range("L43").Value = 32742.25
$ws.Range("N43").Value = -33368.25

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 7593.75
$ws.Range("I99").Value = 7625
$ws.Range("K99").Value = 7625
$ws.Range("M99").Value = -6127
# Row 105
$ws.Range("H105").Value = 11819591
$ws.Range("I105").Value = 715362.9
$ws.Range("J105").Value = 31251990
$ws.Range("K105").Value = 715362.9
$ws.Range("L105").Value = 31251990
$ws.Range("M105").Value = -713615.9
$ws.Range("N105").Value = -31255484
# Row 134
$ws.Range("H134").Value = 1853.174
$ws.Range("I134").Value = 1247.6177
$ws.Range("J134").Value = 3568.9167
$ws.Range("K134").Value = 3742.8531
$ws.Range("L134").Value = 10706.7501
$ws.Range("M134").Value = -1207.8531
$ws.Range("N134").Value = -15776.7501

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 14998.5
$ws.Range("J15").Value = 14998.5
$ws.Range("L15").Value = 14998.5
$ws.Range("N15").Value = -15338.5
# Row 54
$ws.Range("H54").Value = 12116.333
$ws.Range("I54").Value = 12116.333
$ws.Range("K54").Value = 12116.333
$ws.Range("M54").Value = -11458.333
# Row 58
$ws.Range("H58").Value = 2498.6667
$ws.Range("I58").Value = 1907.8
$ws.Range("K58").Value = 1907.8
$ws.Range("M58").Value = -1704.8
# Row 105
$ws.Range("H105").Value = 2506.2
$ws.Range("I105").Value = 2340
$ws.Range("K105").Value = 2340
$ws.Range("M105").Value = -593
# Row 122
$ws.Range("H122").Value = 2090.2307
$ws.Range("I122").Value = 2182.6
$ws.Range("J122").Value = 1782.3334
$ws.Range("K122").Value = 6547.799999999999
$ws.Range("L122").Value = 5347.0002
$ws.Range("M122").Value = -4097.799999999999
$ws.Range("N122").Value = -10247.0002
# Row 132
$ws.Range("H132").Value = 13894356
$ws.Range("I132").Value = 4630
$ws.Range("K132").Value = 13890
$ws.Range("M132").Value = -11360
# Row 136
$ws.Range("H136").Value = 2498.6667
$ws.Range("I136").Value = 1907.8
$ws.Range("K136").Value = 5723.4
$ws.Range("M136").Value = -3173.4

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1048.5834
$ws.Range("I107").Value = 769.25
$ws.Range("J107").Value = 1188.25
$ws.Range("K107").Value = 2307.75
$ws.Range("L107").Value = 3564.75
$ws.Range("M107").Value = -387.75
$ws.Range("N107").Value = -7404.75
# Row 113
$ws.Range("H113").Value = 1699.25
$ws.Range("J113").Value = 1966
$ws.Range("L113").Value = 5898
$ws.Range("N113").Value = -10238

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 27
$ws.Range("H27").Value = 10000000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 122
$ws.Range("H122").Value = 8551899
$ws.Range("I122").Value = 15388414
$ws.Range("K122").Value = 46165242
$ws.Range("M122").Value = -46162792
# Row 126
$ws.Range("H126").Value = 10600.889
$ws.Range("I126").Value = 2032.3334
$ws.Range("K126").Value = 6097.0002
$ws.Range("M126").Value = -3627.0002

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2414.1428
$ws.Range("I7").Value = 1979.8
$ws.Range("K7").Value = 1979.8
$ws.Range("M7").Value = -1867.8
# Row 126
$ws.Range("H126").Value = 2414.1428
$ws.Range("I126").Value = 1979.8
$ws.Range("K126").Value = 5939.4
$ws.Range("M126").Value = -3469.4
# Row 136
$ws.Range("H136").Value = 5167.8696
$ws.Range("I136").Value = 4126.909
$ws.Range("K136").Value = 12380.727
$ws.Range("M136").Value = -9830.726999999999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 8994.5
$ws.Range("J62").Value = 8989
$ws.Range("L62").Value = 8989
$ws.Range("N62").Value = -10237
# Row 65
$ws.Range("H65").Value = 8994.5
$ws.Range("J65").Value = 8989
$ws.Range("L65").Value = 44945
$ws.Range("N65").Value = -51185
# Row 100
$ws.Range("H100").Value = 83334470
$ws.Range("I100").Value = 1313.8334
$ws.Range("K100").Value = 2627.6668
$ws.Range("M100").Value = -2086.6668
# Row 122
$ws.Range("H122").Value = 8336412.5
$ws.Range("I122").Value = 3264.5
$ws.Range("K122").Value = 9793.5
$ws.Range("M122").Value = -7343.5
# Row 132
$ws.Range("H132").Value = 2204.75
$ws.Range("I132").Value = 1905.7778
$ws.Range("K132").Value = 5717.3334
$ws.Range("M132").Value = -3187.3334
# Row 136
$ws.Range("H136").Value = 513133.88
$ws.Range("I136").Value = 2831.95
$ws.Range("K136").Value = 8495.849999999999
$ws.Range("M136").Value = -5945.849999999999
